$wb = $excel.ActiveWorkbook

# Insert the new "unit test checklist" sheet between "First year" and "Second Year"
$firstYear = $wb.Worksheets.Item("First year")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $firstYear)
$newSheet.Name = "unit test checklist"

# Column A values (shared-string order must match: Linked list, create list, ..., clear)
$newSheet.Range("A1").Value = "Linked list"
$newSheet.Range("A2").Value = "create list"
$newSheet.Range("A3").Value = "pushback"
$newSheet.Range("A4").Value = "push front"
$newSheet.Range("A5").Value = "insert"
$newSheet.Range("A6").Value = "begin"
$newSheet.Range("A7").Value = "end"
$newSheet.Range("A8").Value = "first"
$newSheet.Range("A9").Value = "last"
$newSheet.Range("A10").Value = "count"
$newSheet.Range("A11").Value = "erase"
$newSheet.Range("A12").Value = "remove"
$newSheet.Range("A13").Value = "popBack"
$newSheet.Range("A14").Value = "popFront"
$newSheet.Range("A15").Value = "empty"
$newSheet.Range("A16").Value = "clear"

# Column headers B1/C1 ("Binary Tree" is written last in the source so it
# lands at the end of the shared-string table)
$newSheet.Range("B1").Value = "Binary Tree"
$newSheet.Range("C1").Value = "Hash Map"

# First year sheet's selection moved to the "Production" section header
$ws1 = $wb.Worksheets.Item("First year")
$ws1.Range("A19:F19").Select()

# First year sheet column E got a bit wider
$ws1.Columns("E").ColumnWidth = 78.85546875

# Activate the new sheet last so it is the visible/selected tab
$newSheet.Activate()
